$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 3 (Blink Ops / Peter Draper entry) and shift the remaining rows up,
# matching the weekly "last 30 days" report refresh.
$ws.Rows("3").Delete()

# Two rows (Dash0 / Blink Ops postings) kept their existing Status while the
# other columns moved up with the shift; restore Status and set the refreshed
# Action Date for those two rows to match the latest export.
$ws.Range("E9").Value = "1st Interview"
$ws.Range("F9").Value = 45993
$ws.Range("E12").Value = "1st Interview"
$ws.Range("F12").Value = 45993
